$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.206.81'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.34%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.602.45'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.11%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.001'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '303.11'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.65%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3781'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '51.81'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.43%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3616'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.76%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.266'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.69%  '

$ws.Range("E11").Value = '  -0.05%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08127'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.13%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.60'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.61%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.587'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.13%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.398'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.47%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001249'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.78%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.601.56'
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '93.81'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.31%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06880'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.23%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.06'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.536'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.37%  '

$ws.Range("E22").Value = '  +0.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.97'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.37%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.202.02'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.29%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.389'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.82%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.975'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +9.69%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.19'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.49%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '149.53'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.32%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.254'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.81%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.82'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.53%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.401'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.27%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.884'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.77%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.780.11'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.05%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9754'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.77%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07535'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.79%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.30'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.65%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02720'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.52%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.128'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.75%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2507'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.65%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08802'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.08%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7116'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.31%  '

$ws.Range("E42").Value = '  -1.79%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.52'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.92%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.37'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.18%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6547'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.87%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.309'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.38%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.014'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.88%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.36'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.50%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07957'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.16%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.204'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.48%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.220'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.05%  '
